$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2256.8667
$ws.Range("I15").Value = 2256.8667
$ws.Range("K15").Value = 6770.6001
$ws.Range("M15").Value = -6601.6001
$ws.Range("H17").Value = 1875
$ws.Range("J17").Value = 1875
$ws.Range("L17").Value = 5625
$ws.Range("N17").Value = -5961
$ws.Range("H40").Value = 4055.4443
$ws.Range("I40").Value = 3739.8
$ws.Range("J40").Value = 4450
$ws.Range("K40").Value = 3739.8
$ws.Range("L40").Value = 4450
$ws.Range("M40").Value = -3564.8
$ws.Range("N40").Value = -4800
$ws.Range("H64").Value = 4043.5
$ws.Range("I64").Value = 3498.5
$ws.Range("J64").Value = 4861
$ws.Range("K64").Value = 3498.5
$ws.Range("L64").Value = 4861
$ws.Range("M64").Value = -3250.5
$ws.Range("N64").Value = -5357
$ws.Range("H67").Value = 4043.5
$ws.Range("I67").Value = 3498.5
$ws.Range("J67").Value = 4861
$ws.Range("K67").Value = 3498.5
$ws.Range("L67").Value = 4861
$ws.Range("M67").Value = -2640.5
$ws.Range("N67").Value = -6577
$ws.Range("H69").Value = 17000.572
$ws.Range("J69").Value = 18335
$ws.Range("L69").Value = 55005
$ws.Range("N69").Value = -56753
$ws.Range("H72").Value = 17000.572
$ws.Range("J72").Value = 18335
$ws.Range("L72").Value = 165015
$ws.Range("N72").Value = -173751
$ws.Range("H76").Value = 4205.0386
$ws.Range("I76").Value = 3365.611
$ws.Range("K76").Value = 3365.611
$ws.Range("M76").Value = -3050.611
$ws.Range("H79").Value = 4205.0386
$ws.Range("I79").Value = 3365.611
$ws.Range("K79").Value = 3365.611
$ws.Range("M79").Value = -2273.611
$ws.Range("H100").Value = 5269.467
$ws.Range("I100").Value = 1775.5
$ws.Range("K100").Value = 1775.5
$ws.Range("M100").Value = -1234.5
$ws.Range("H112").Value = 2155.3076
$ws.Range("J112").Value = 2243.3333
$ws.Range("L112").Value = 6729.999899999999
$ws.Range("N112").Value = -8945.999899999999
$ws.Range("H132").Value = 2202.5217
$ws.Range("I132").Value = 1426.9412
$ws.Range("K132").Value = 4280.8236
$ws.Range("M132").Value = -1750.8236
$ws.Range("H138").Value = 2514
$ws.Range("I138").Value = 1017.9231
$ws.Range("K138").Value = 3053.7693
$ws.Range("M138").Value = 2086.2307
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 16131709
$ws.Range("I45").Value = 26317696
$ws.Range("K45").Value = 26317696
$ws.Range("M45").Value = -26317319
$ws.Range("H132").Value = 3541.0312
$ws.Range("I132").Value = 3110.5518
$ws.Range("J132").Value = 7702.3335
$ws.Range("K132").Value = 9331.6554
$ws.Range("L132").Value = 23107.0005
$ws.Range("M132").Value = -6801.6554
$ws.Range("N132").Value = -28167.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19971
$ws.Range("I26").Value = 19971
$ws.Range("K26").Value = 19971
$ws.Range("M26").Value = -19679
$ws.Range("H134").Value = 3849516.2
$ws.Range("I134").Value = 1681.4667
$ws.Range("K134").Value = 5044.4001
$ws.Range("M134").Value = -2509.4001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 80377.8
$ws.Range("J118").Value = 80377.8
$ws.Range("L118").Value = 80377.8
$ws.Range("N118").Value = -83691.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 7014.7144
$ws.Range("I16").Value = 7500
$ws.Range("J16").Value = 6933.8335
$ws.Range("K16").Value = 22500
$ws.Range("L16").Value = 20801.5005
$ws.Range("M16").Value = -22327
$ws.Range("N16").Value = -21147.5005
$ws.Range("H34").Value = 528.3333
$ws.Range("I34").Value = 517.5
$ws.Range("J34").Value = 550
$ws.Range("K34").Value = 1552.5
$ws.Range("L34").Value = 1650
$ws.Range("M34").Value = -1468.5
$ws.Range("N34").Value = -1818
$ws.Range("H56").Value = 6027.8887
$ws.Range("I56").Value = 6027.8887
$ws.Range("K56").Value = 6027.8887
$ws.Range("M56").Value = -5497.8887
$ws.Range("H113").Value = 1034
$ws.Range("J113").Value = 1129.1538
$ws.Range("L113").Value = 3387.4614
$ws.Range("N113").Value = -7727.4614
$ws.Range("H117").Value = 1880.5
$ws.Range("J117").Value = 2027.7142
$ws.Range("L117").Value = 6083.142599999999
$ws.Range("N117").Value = -12967.1426
$ws.Range("H129").Value = 25645196
$ws.Range("I129").Value = 2600
$ws.Range("J129").Value = 33337976
$ws.Range("K129").Value = 7800
$ws.Range("L129").Value = 100013928
$ws.Range("M129").Value = -2800
$ws.Range("N129").Value = -100023928
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 82501.875
$ws.Range("I35").Value = 97753.75
$ws.Range("K35").Value = 97753.75
$ws.Range("M35").Value = -97455.75
$ws.Range("H39").Value = 52500
$ws.Range("J39").Value = 52500
$ws.Range("L39").Value = 52500
$ws.Range("N39").Value = -53564
$ws.Range("H59").Value = 21755
$ws.Range("I59").Value = 7019
$ws.Range("J59").Value = 26667
$ws.Range("K59").Value = 7019
$ws.Range("L59").Value = 26667
$ws.Range("M59").Value = -6436
$ws.Range("N59").Value = -27833
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 14000
$ws.Range("J13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("N13").Value = -10280
$ws.Range("H16").Value = 1843.1818
$ws.Range("I16").Value = 1877.5
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1877.5
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1707.5
$ws.Range("N16").Value = -1840
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 706.1875
$ws.Range("I14").Value = 200
$ws.Range("J14").Value = 2224.75
$ws.Range("K14").Value = 200
$ws.Range("L14").Value = 2224.75
$ws.Range("M14").Value = -32
$ws.Range("N14").Value = -2560.75
$ws.Range("H38").Value = 14814
$ws.Range("I38").Value = 12418.667
$ws.Range("J38").Value = 22000
$ws.Range("K38").Value = 12418.667
$ws.Range("L38").Value = 22000
$ws.Range("M38").Value = -11945.667
$ws.Range("N38").Value = -22946
$ws.Range("H49").Value = 31245
$ws.Range("J49").Value = 31245
$ws.Range("L49").Value = 31245
$ws.Range("N49").Value = -31705
$ws.Range("H82").Value = 34494.5
$ws.Range("J82").Value = 36329.332
$ws.Range("L82").Value = 36329.332
$ws.Range("N82").Value = -37095.332
$ws.Range("H85").Value = 34494.5
$ws.Range("J85").Value = 36329.332
$ws.Range("L85").Value = 36329.332
$ws.Range("N85").Value = -38981.332
$ws.Range("H132").Value = 297335.4
$ws.Range("I132").Value = 3080
$ws.Range("K132").Value = 9240
$ws.Range("M132").Value = -6710
